$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Persona -> Organización
$ws.Range("A2").Value = "Organización"
$ws.Range("B2").Value = "Atributos de una organización a la cual o donde se realiza la denuncia."

# Row 3: Organización -> Denuncia
$ws.Range("A3").Value = "Denuncia"
$ws.Range("B3").Value = "Atributos de la denuncia realizada por la persona en la organización."

# Remove old row 4 (now redundant since its content moved to row 3)
$ws.Rows.Item(4).Delete()
